$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write strings in the desired "first seen" order so the shared-string
# table comes out as Carrots(0), Potatos(1), Tomatos(2), Corn(3).
$ws.Range("A3").Value = "Carrots"
$ws.Range("B3").Value = 300

$ws.Range("A2").Value = "Potatos"
$ws.Range("B2").Value = 400

$ws.Range("A1").Value = "Tomatos"
$ws.Range("B1").Value = 500

$ws.Range("A4").Value = "Corn"
$ws.Range("B4").Value = 200

# The amount cells in rows 1-2 lose their inherited style once retyped.
$null = $ws.Range("B1").ClearFormats()
$null = $ws.Range("B2").ClearFormats()

# Leave the selection where the user last clicked.
$null = $ws.Range("E4").Select()
